# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-parsed as a number by
# Excel (single decimal point) are pre-formatted as Text so they stay strings,
# matching the original inline-string cell type.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "58.123.00"
$ws.Range("E2").Value = "  -1.91%  "

$ws.Range("D3").Value = "2.468.90"
$ws.Range("E3").Value = "  -2.27%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "518.82"
$ws.Range("E5").Value = "  -3.41%  "

$ws.Range("D6").Value = "132.34"
$ws.Range("E6").Value = "  -4.26%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -1.78%  "

$ws.Range("E10").Value = "  -0.65%  "

$ws.Range("D11").Value = "5.38"
$ws.Range("E11").Value = "  +0.40%  "

$ws.Range("E12").Value = "  -2.18%  "

$ws.Range("D13").Value = "2.908.89"
$ws.Range("E13").Value = "  -2.31%  "

$ws.Range("D14").Value = "58.099.73"
$ws.Range("E14").Value = "  -1.84%  "

$ws.Range("D15").Value = "22.09"
$ws.Range("E15").Value = "  -4.69%  "

$ws.Range("E16").Value = "  -2.77%  "

$ws.Range("D17").Value = "2.472.96"
$ws.Range("E17").Value = "  -2.18%  "

$ws.Range("E18").Value = "  -2.48%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "320.43"
$ws.Range("E19").Value = "  -1.64%  "

$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "4.18"
$ws.Range("E20").Value = "  -2.75%  "

$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("D22").Value = "5.75"
$ws.Range("E22").Value = "  -3.81%  "

$ws.Range("D23").Value = "64.35"
$ws.Range("E23").Value = "  -1.92%  "

$ws.Range("E24").Value = "  -3.71%  "

$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.30%  "

$ws.Range("D26").Value = "0.162"
$ws.Range("E26").Value = "  -3.53%  "

$ws.Range("E27").Value = "  -3.77%  "

$ws.Range("E28").Value = "  -3.78%  "

$ws.Range("D29").Value = "6.39"
$ws.Range("E29").Value = "  -5.54%  "

$ws.Range("E30").Value = "  -4.90%  "

$ws.Range("D31").Value = "165.52"
$ws.Range("E31").Value = "  +0.77%  "

$ws.Range("E32").Value = "  -4.34%  "

$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("E34").Value = "  -0.11%  "

$ws.Range("D35").Value = "18.14"
$ws.Range("E35").Value = "  -1.97%  "

$ws.Range("E36").Value = "  -9.66%  "

$ws.Range("E37").Value = "  -3.49%  "

$ws.Range("E38").Value = "  -4.16%  "

$ws.Range("D39").Value = "0.796"
$ws.Range("E39").Value = "  -2.90%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "276.02"
$ws.Range("E40").Value = "  -4.26%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "3.47"

$ws.Range("D42").Value = "5.01"
$ws.Range("E42").Value = "  -4.26%  "

$ws.Range("D43").Value = "0.593"
$ws.Range("E43").Value = "  -3.19%  "

$ws.Range("D44").Value = "126.51"
$ws.Range("E44").Value = "  -4.68%  "

$ws.Range("E45").Value = "  -2.81%  "

$ws.Range("E46").Value = "  -3.69%  "

$ws.Range("E47").Value = "  -3.55%  "

$ws.Range("D48").Value = "17.15"
$ws.Range("E48").Value = "  -1.69%  "

$ws.Range("D49").Value = "1.733.12"
$ws.Range("E49").Value = "  -1.74%  "

$ws.Range("D51").Value = "4.66"
$ws.Range("E51").Value = "  -2.13%  "
